# "Fixed LMS end Login And Logout"
# - Updates the login credentials on the TestCredentials sheet:
#     * LearnerLogin  -> email becomes a mailto hyperlink + new password
#     * CALogin       -> credentials cleared out
#     * ManagerLogin  -> email becomes a mailto hyperlink + new password
# - TestCredentials becomes the active/selected sheet (was CreateAssessment).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TestCredentials")

# --- Row 2 (LearnerLogin): new e-mail (as hyperlink) + new password ---
$ws1.Range("C2").Value = "y.adityaprasad@gmail.com"
$ws1.Hyperlinks.Add($ws1.Range("C2"), "mailto:y.adityaprasad@gmail.com")
$ws1.Range("D2").Value = "Password@123"
$ws1.Range("D2").Borders.LineStyle = -4142

# --- Row 3 (CALogin): credentials removed ---
$ws1.Range("C3").ClearContents()
$ws1.Range("D3").ClearContents()

# --- Row 4 (ManagerLogin): new e-mail (as hyperlink) + new password ---
$ws1.Range("C4").Value = "khmanagertest@khcodelab.com"
$ws1.Hyperlinks.Add($ws1.Range("C4"), "mailto:khmanagertest@khcodelab.com")
$ws1.Range("D4").Value = "Password@123"
$ws1.Range("D4").Borders.LineStyle = -4142

# Move the sheet selection / view to C4 on this sheet.
$ws1.Range("C4").Select()

# TestCredentials is now the active tab (previously CreateAssessment).
$ws1.Activate()

Write-Output "edit complete"
